$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-03-28 Friday", $true, $true, $false, $false, $false, $true, 1, $false, "2025-03-29 Saturday", 2) | Out-Null
$d.Content.Find.Execute("70+10=", $true, $true, $false, $false, $false, $true, 1, $false, "3+25=", 2) | Out-Null
$d.Content.Find.Execute("74-32=", $true, $true, $false, $false, $false, $true, 1, $false, "44-38=", 2) | Out-Null
$d.Content.Find.Execute("73-21=", $true, $true, $false, $false, $false, $true, 1, $false, "69+16=", 2) | Out-Null
$d.Content.Find.Execute("72-34=", $true, $true, $false, $false, $false, $true, 1, $false, "40-14=", 2) | Out-Null
$d.Content.Find.Execute("93-45=", $true, $true, $false, $false, $false, $true, 1, $false, "26-2=", 2) | Out-Null
$d.Content.Find.Execute("42-26=", $true, $true, $false, $false, $false, $true, 1, $false, "19+2=", 2) | Out-Null
$d.Content.Find.Execute("70-41=", $true, $true, $false, $false, $false, $true, 1, $false, "55-27=", 2) | Out-Null
$d.Content.Find.Execute("13+27=", $true, $true, $false, $false, $false, $true, 1, $false, "93-9=", 2) | Out-Null
$d.Content.Find.Execute("69-43=", $true, $true, $false, $false, $false, $true, 1, $false, "37+25=", 2) | Out-Null
$d.Content.Find.Execute("98-75=", $true, $true, $false, $false, $false, $true, 1, $false, "92-88=", 2) | Out-Null
$d.Content.Find.Execute("12+32=", $true, $true, $false, $false, $false, $true, 1, $false, "46+0=", 2) | Out-Null
$d.Content.Find.Execute("35+49=", $true, $true, $false, $false, $false, $true, 1, $false, "60-22=", 2) | Out-Null
$d.Content.Find.Execute("60+7=", $true, $true, $false, $false, $false, $true, 1, $false, "91-27=", 2) | Out-Null
$d.Content.Find.Execute("72-59=", $true, $true, $false, $false, $false, $true, 1, $false, "54-21=", 2) | Out-Null
$d.Content.Find.Execute("89-82=", $true, $true, $false, $false, $false, $true, 1, $false, "90-63=", 2) | Out-Null
$d.Content.Find.Execute("73+4=", $true, $true, $false, $false, $false, $true, 1, $false, "34+30=", 2) | Out-Null
$d.Content.Find.Execute("99-19=", $true, $true, $false, $false, $false, $true, 1, $false, "24+54=", 2) | Out-Null
$d.Content.Find.Execute("45+52=", $true, $true, $false, $false, $false, $true, 1, $false, "32+30=", 2) | Out-Null
$d.Content.Find.Execute("97-71=", $true, $true, $false, $false, $false, $true, 1, $false, "70+5=", 2) | Out-Null
$d.Content.Find.Execute("46-14=", $true, $true, $false, $false, $false, $true, 1, $false, "45-15=", 2) | Out-Null
$d.Content.Find.Execute("61-35=", $true, $true, $false, $false, $false, $true, 1, $false, "57-47=", 2) | Out-Null
$d.Content.Find.Execute("45+8=", $true, $true, $false, $false, $false, $true, 1, $false, "96-38=", 2) | Out-Null
$d.Content.Find.Execute("66-41=", $true, $true, $false, $false, $false, $true, 1, $false, "28+57=", 2) | Out-Null
$d.Content.Find.Execute("65-62=", $true, $true, $false, $false, $false, $true, 1, $false, "28+34=", 2) | Out-Null
$d.Content.Find.Execute("26+13=", $true, $true, $false, $false, $false, $true, 1, $false, "83-14=", 2) | Out-Null
$d.Content.Find.Execute("12+19=", $true, $true, $false, $false, $false, $true, 1, $false, "53+42=", 2) | Out-Null
$d.Content.Find.Execute("49-18=", $true, $true, $false, $false, $false, $true, 1, $false, "74-58=", 2) | Out-Null
$d.Content.Find.Execute("99-27=", $true, $true, $false, $false, $false, $true, 1, $false, "3+38=", 2) | Out-Null
$d.Content.Find.Execute("54+34=", $true, $true, $false, $false, $false, $true, 1, $false, "52-40=", 2) | Out-Null
$d.Content.Find.Execute("1+10=", $true, $true, $false, $false, $false, $true, 1, $false, "28-24=", 2) | Out-Null
$d.Content.Find.Execute("7+69=", $true, $true, $false, $false, $false, $true, 1, $false, "68+29=", 2) | Out-Null
$d.Content.Find.Execute("30-30=", $true, $true, $false, $false, $false, $true, 1, $false, "63+16=", 2) | Out-Null
$d.Content.Find.Execute("22+70=", $true, $true, $false, $false, $false, $true, 1, $false, "57-43=", 2) | Out-Null
$d.Content.Find.Execute("22+40=", $true, $true, $false, $false, $false, $true, 1, $false, "70-13=", 2) | Out-Null
$d.Content.Find.Execute("34-5=", $true, $true, $false, $false, $false, $true, 1, $false, "77-47=", 2) | Out-Null
$d.Content.Find.Execute("42-15=", $true, $true, $false, $false, $false, $true, 1, $false, "14+71=", 2) | Out-Null
$d.Content.Find.Execute("6+80=", $true, $true, $false, $false, $false, $true, 1, $false, "90-71=", 2) | Out-Null
$d.Content.Find.Execute("11+25=", $true, $true, $false, $false, $false, $true, 1, $false, "28-15=", 2) | Out-Null
$d.Content.Find.Execute("33-20=", $true, $true, $false, $false, $false, $true, 1, $false, "10+86=", 2) | Out-Null
$d.Content.Find.Execute("27-24=", $true, $true, $false, $false, $false, $true, 1, $false, "37+13=", 2) | Out-Null
$d.Content.Find.Execute("54-33=", $true, $true, $false, $false, $false, $true, 1, $false, "6+4=", 2) | Out-Null
$d.Content.Find.Execute("79+11=", $true, $true, $false, $false, $false, $true, 1, $false, "33+35=", 2) | Out-Null
$d.Content.Find.Execute("89-79=", $true, $true, $false, $false, $false, $true, 1, $false, "35+27=", 2) | Out-Null
$d.Content.Find.Execute("54-40=", $true, $true, $false, $false, $false, $true, 1, $false, "16+7=", 2) | Out-Null
$d.Content.Find.Execute("80-69=", $true, $true, $false, $false, $false, $true, 1, $false, "27-20=", 2) | Out-Null
$d.Content.Find.Execute("77-51=", $true, $true, $false, $false, $false, $true, 1, $false, "38+20=", 2) | Out-Null
$d.Content.Find.Execute("6+33=", $true, $true, $false, $false, $false, $true, 1, $false, "97-67=", 2) | Out-Null
$d.Content.Find.Execute("39+24=", $true, $true, $false, $false, $false, $true, 1, $false, "19+10=", 2) | Out-Null
$d.Content.Find.Execute("35+63=", $true, $true, $false, $false, $false, $true, 1, $false, "61+2=", 2) | Out-Null
$d.Content.Find.Execute("52-30=", $true, $true, $false, $false, $false, $true, 1, $false, "31+45=", 2) | Out-Null
$d.Content.Find.Execute("1+90=", $true, $true, $false, $false, $false, $true, 1, $false, "6-2=", 2) | Out-Null
$d.Content.Find.Execute("97-9=", $true, $true, $false, $false, $false, $true, 1, $false, "35-13=", 2) | Out-Null
$d.Content.Find.Execute("88-18=", $true, $true, $false, $false, $false, $true, 1, $false, "33+31=", 2) | Out-Null
$d.Content.Find.Execute("65-22=", $true, $true, $false, $false, $false, $true, 1, $false, "52+6=", 2) | Out-Null
$d.Content.Find.Execute("48+10=", $true, $true, $false, $false, $false, $true, 1, $false, "27-6=", 2) | Out-Null
$d.Content.Find.Execute("65-43=", $true, $true, $false, $false, $false, $true, 1, $false, "96-19=", 2) | Out-Null
$d.Content.Find.Execute("16-8=", $true, $true, $false, $false, $false, $true, 1, $false, "83-12=", 2) | Out-Null
$d.Content.Find.Execute("56+21=", $true, $true, $false, $false, $false, $true, 1, $false, "85-33=", 2) | Out-Null
$d.Content.Find.Execute("9+45=", $true, $true, $false, $false, $false, $true, 1, $false, "86-11=", 2) | Out-Null
$d.Content.Find.Execute("91-86=", $true, $true, $false, $false, $false, $true, 1, $false, "38+45=", 2) | Out-Null
$d.Content.Find.Execute("74+13=", $true, $true, $false, $false, $false, $true, 1, $false, "59+34=", 2) | Out-Null
$d.Content.Find.Execute("70-68=", $true, $true, $false, $false, $false, $true, 1, $false, "57-6=", 2) | Out-Null
$d.Content.Find.Execute("85-31=", $true, $true, $false, $false, $false, $true, 1, $false, "39-28=", 2) | Out-Null
$d.Content.Find.Execute("66+13=", $true, $true, $false, $false, $false, $true, 1, $false, "48-45=", 2) | Out-Null
$d.Content.Find.Execute("25+22=", $true, $true, $false, $false, $false, $true, 1, $false, "77-21=", 2) | Out-Null
$d.Content.Find.Execute("24+25=", $true, $true, $false, $false, $false, $true, 1, $false, "90+0=", 2) | Out-Null
$d.Content.Find.Execute("53-14=", $true, $true, $false, $false, $false, $true, 1, $false, "31+45=", 2) | Out-Null
$d.Content.Find.Execute("63+19=", $true, $true, $false, $false, $false, $true, 1, $false, "31+22=", 2) | Out-Null
$d.Content.Find.Execute("38+4=", $true, $true, $false, $false, $false, $true, 1, $false, "26-16=", 2) | Out-Null
$d.Content.Find.Execute("22+22=", $true, $true, $false, $false, $false, $true, 1, $false, "25+65=", 2) | Out-Null
$d.Content.Find.Execute("64-56=", $true, $true, $false, $false, $false, $true, 1, $false, "78+2=", 2) | Out-Null
$d.Content.Find.Execute("12+84=", $true, $true, $false, $false, $false, $true, 1, $false, "10+20=", 2) | Out-Null
$d.Content.Find.Execute("43+32=", $true, $true, $false, $false, $false, $true, 1, $false, "39+32=", 2) | Out-Null
$d.Content.Find.Execute("25+46=", $true, $true, $false, $false, $false, $true, 1, $false, "28+27=", 2) | Out-Null
$d.Content.Find.Execute("14+45=", $true, $true, $false, $false, $false, $true, 1, $false, "87-52=", 2) | Out-Null
$d.Content.Find.Execute("20+32=", $true, $true, $false, $false, $false, $true, 1, $false, "50+28=", 2) | Out-Null
$d.Content.Find.Execute("76-69=", $true, $true, $false, $false, $false, $true, 1, $false, "60+34=", 2) | Out-Null
$d.Content.Find.Execute("18+73=", $true, $true, $false, $false, $false, $true, 1, $false, "15+8=", 2) | Out-Null
$d.Content.Find.Execute("27+37=", $true, $true, $false, $false, $false, $true, 1, $false, "78-62=", 2) | Out-Null
$d.Content.Find.Execute("10+39=", $true, $true, $false, $false, $false, $true, 1, $false, "81-9=", 2) | Out-Null
$d.Content.Find.Execute("85+12=", $true, $true, $false, $false, $false, $true, 1, $false, "32+4=", 2) | Out-Null
$d.Content.Find.Execute("94-29=", $true, $true, $false, $false, $false, $true, 1, $false, "26-17=", 2) | Out-Null
$d.Content.Find.Execute("1-0=", $true, $true, $false, $false, $false, $true, 1, $false, "66-64=", 2) | Out-Null
$d.Content.Find.Execute("92+2=", $true, $true, $false, $false, $false, $true, 1, $false, "40-35=", 2) | Out-Null
$d.Content.Find.Execute("46-8=", $true, $true, $false, $false, $false, $true, 1, $false, "0+1=", 2) | Out-Null
$d.Content.Find.Execute("72+23=", $true, $true, $false, $false, $false, $true, 1, $false, "51-43=", 2) | Out-Null
$d.Content.Find.Execute("83-5=", $true, $true, $false, $false, $false, $true, 1, $false, "80-56=", 2) | Out-Null
$d.Content.Find.Execute("7+30=", $true, $true, $false, $false, $false, $true, 1, $false, "21-15=", 2) | Out-Null
$d.Content.Find.Execute("74-37=", $true, $true, $false, $false, $false, $true, 1, $false, "58-13=", 2) | Out-Null
$d.Content.Find.Execute("9+81=", $true, $true, $false, $false, $false, $true, 1, $false, "38-24=", 2) | Out-Null
$d.Content.Find.Execute("50+25=", $true, $true, $false, $false, $false, $true, 1, $false, "36+46=", 2) | Out-Null
$d.Content.Find.Execute("9+66=", $true, $true, $false, $false, $false, $true, 1, $false, "79-78=", 2) | Out-Null
$d.Content.Find.Execute("30+1=", $true, $true, $false, $false, $false, $true, 1, $false, "39+54=", 2) | Out-Null
$d.Content.Find.Execute("99+0=", $true, $true, $false, $false, $false, $true, 1, $false, "33+27=", 2) | Out-Null
$d.Content.Find.Execute("57+42=", $true, $true, $false, $false, $false, $true, 1, $false, "20-3=", 2) | Out-Null
$d.Content.Find.Execute("90-50=", $true, $true, $false, $false, $false, $true, 1, $false, "43-40=", 2) | Out-Null
$d.Content.Find.Execute("51-26=", $true, $true, $false, $false, $false, $true, 1, $false, "55+25=", 2) | Out-Null
$d.Content.Find.Execute("49-33=", $true, $true, $false, $false, $false, $true, 1, $false, "60+8=", 2) | Out-Null
$d.Content.Find.Execute("79-31=", $true, $true, $false, $false, $false, $true, 1, $false, "42+55=", 2) | Out-Null
$d.Content.Find.Execute("40+41=", $true, $true, $false, $false, $false, $true, 1, $false, "83-6=", 2) | Out-Null
